# Daily update at 8 AM UTC
# Appends the next day's row (row 28) to the "Wins Over Time" tracking
# sheet, and restores the now-second-to-last row (27) to the regular
# date/time number format, moving the "last row" date-only formatting
# down to the newly appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 is no longer the last row, so give it back the same number
# format used by the rest of the table (copy it from row 26, which is
# already in the regular format).
$ws.Range("A27").NumberFormat = $ws.Range("A26").NumberFormat

# Append the new day's data as row 28.
$ws.Range("A28").Value = 45768
$ws.Range("B28").Value = 111
$ws.Range("C28").Value = 114
$ws.Range("D28").Value = 112

# The new last row gets the date-only number format.
$ws.Range("A28").NumberFormat = "YYYY-MM-DD"
